$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.418.66"
$ws.Range("E2").Value = "  +3.96%  "
$ws.Range("D3").Value = "1.588.65"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  -0.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.90%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").Value = "1.814.68"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "1.584.00"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("E15").Value = "  +2.52%  "
$ws.Range("D16").Value = "28.412.56"
$ws.Range("E16").Value = "  +4.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.62%  "
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "1.411.77"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("E36").Value = "  -4.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  +8.06%  "
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.818"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("E42").Value = "  -2.86%  "
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("E45").Value = "  +5.60%  "
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "1.727.50"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("E49").Value = "  +4.77%  "
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +15.56%  "
